$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Estadisticos 1P" - row 7 (TEMAS DE BIOLOGÍA CONTEMPORÁNEA / 6ALCV)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Cells.Item(7,4).Value = 3        # D7: Blancos     13 -> 3
$ws1.Cells.Item(7,6).Value = 22       # F7: Aprobados   12 -> 22
$ws1.Cells.Item(7,7).Value = 88       # G7: Por_Apro    48 -> 88
$ws1.Cells.Item(7,8).Value = 7.8      # H7: Promedio     8 -> 7.8

# ---------------------------------------------------------------------------
# 2) "Estadisticos 2P" - row 7
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Cells.Item(7,4).Value = 3        # D7: Blancos     25 -> 3
$ws2.Cells.Item(7,5).Value = 0        # E7: Reprobados  12 -> 0
$ws2.Cells.Item(7,6).Value = 22       # F7: Aprobados    0 -> 22
$ws2.Cells.Item(7,7).Value = 88       # G7: Por_Apro     0 -> 88
$ws2.Cells.Item(7,8).Value = 7.8      # H7: Promedio   (new) 7.8

# ---------------------------------------------------------------------------
# 3) "Estadisticos Final" - row 7
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Cells.Item(7,4).Value = 3        # D7: Blancos     13 -> 3
$ws3.Cells.Item(7,6).Value = 22       # F7: Aprobados   12 -> 22
$ws3.Cells.Item(7,7).Value = 88       # G7: Por_Apro    48 -> 88
# H7 (Promedio) stays 8 - unchanged

# ---------------------------------------------------------------------------
# 4) "Rescatables" - table shrinks from 20 data rows to 4 data rows, and the
#    remaining rows get new people / values.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")

# Drop everything below the new last row (A6:G20) so the used range shrinks
# to A1:G5, matching the new dimension.
$ws4.Range("A6:G20").ClearContents()

# Row 2
$ws4.Cells.Item(2,1).Value = 20330051920137
$ws4.Cells.Item(2,2).Value = "MAYAHUA"
$ws4.Cells.Item(2,3).Value = "XOCHIQUISQUI"
$ws4.Cells.Item(2,4).Value = "DAMARIS"
$ws4.Cells.Item(2,5).Value = "INGLÉS II"
$ws4.Cells.Item(2,6).Value = "2ARHV"
$ws4.Cells.Item(2,7).Value = 2

# Row 3
$ws4.Cells.Item(3,1).Value = 20330051920325
$ws4.Cells.Item(3,2).Value = "HERNANDEZ"
$ws4.Cells.Item(3,3).Value = "CARRILLO"
$ws4.Cells.Item(3,4).Value = "ANGEL DAVID"
$ws4.Cells.Item(3,5).Value = "INGLÉS II"
$ws4.Cells.Item(3,6).Value = "2ASV"
$ws4.Cells.Item(3,7).Value = 2

# Row 4
$ws4.Cells.Item(4,1).Value = 18330051920245
$ws4.Cells.Item(4,2).Value = "DE JESUS"
$ws4.Cells.Item(4,3).Value = "ISIDRO"
$ws4.Cells.Item(4,4).Value = "MONSERRAT"
$ws4.Cells.Item(4,5).Value = "TEMAS DE BIOLOGÍA CONTEMPORÁNEA"
$ws4.Cells.Item(4,6).Value = "6ALCV"
$ws4.Cells.Item(4,7).Value = 2

# Row 5
$ws4.Cells.Item(5,1).Value = 20330051920061
$ws4.Cells.Item(5,2).Value = "ROBLES"
$ws4.Cells.Item(5,3).Value = "IXMATLAHUA"
$ws4.Cells.Item(5,4).Value = "ALAN URIEL"
$ws4.Cells.Item(5,5).Value = "INGLÉS II"
$ws4.Cells.Item(5,6).Value = "2AEV"
$ws4.Cells.Item(5,7).Value = 1

Write-Output "edits applied"
